$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBF = New-Object "object[,]" 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.05829982881462
$dataBF[0,2] = 1.062805529706471
$dataBF[0,3] = 1.054319396435029
$dataBF[0,4] = 1.070992466017058
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.059565133249015
$dataBF[1,2] = 1.063843583855984
$dataBF[1,3] = 1.055409297003092
$dataBF[1,4] = 1.072223708514826
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.060383472150324
$dataBF[2,2] = 1.064514882631411
$dataBF[2,3] = 1.056114403460226
$dataBF[2,4] = 1.073020403565456
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.060727409255372
$dataBF[3,2] = 1.064797004348009
$dataBF[3,3] = 1.056410800140439
$dataBF[3,4] = 1.073355335811718
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.060785152456982
$dataBF[4,2] = 1.064844368456689
$dataBF[4,3] = 1.056460564728627
$dataBF[4,4] = 1.073411572562233
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.060388068215194
$dataBF[5,2] = 1.064518652716384
$dataBF[5,3] = 1.056118364043415
$dataBF[5,4] = 1.073024878937642
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.058727528300925
$dataBF[6,2] = 1.063156426956001
$dataBF[6,3] = 1.054687761718182
$dataBF[6,4] = 1.071408571589316
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.055798280867718
$dataBF[7,2] = 1.060752944710537
$dataBF[7,3] = 1.0521657762173
$dataBF[7,4] = 1.068560335779567
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.053843154642077
$dataBF[8,2] = 1.05914846583679
$dataBF[8,3] = 1.05048362545457
$dataBF[8,4] = 1.06666132376671
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.052995981265159
$dataBF[9,2] = 1.058453174997814
$dataBF[9,3] = 1.049755016404325
$dataBF[9,4] = 1.065838957660856
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.052681211232783
$dataBF[10,2] = 1.058194829759056
$dataBF[10,3] = 1.049484342428955
$dataBF[10,4] = 1.065533479511561
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.052748734695842
$dataBF[11,2] = 1.058250249489782
$dataBF[11,3] = 1.049542404584044
$dataBF[11,4] = 1.06559900632038
$dataBF[12,0] = 1.02
$dataBF[12,1] = 1.052969964174747
$dataBF[12,2] = 1.058431821809928
$dataBF[12,3] = 1.049732643158231
$dataBF[12,4] = 1.065813707057577
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.053106258650885
$dataBF[13,2] = 1.05854368341586
$dataBF[13,3] = 1.049849850572196
$dataBF[13,4] = 1.065945989235688
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.053899366014892
$dataBF[14,2] = 1.059194598433305
$dataBF[14,3] = 1.050531975884755
$dataBF[14,4] = 1.066715899566037
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.054396700764739
$dataBF[15,2] = 1.059602754003414
$dataBF[15,3] = 1.050959792920843
$dataBF[15,4] = 1.067198820401246
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.054686730932832
$dataBF[16,2] = 1.05984077218973
$dataBF[16,3] = 1.051209309631522
$dataBF[16,4] = 1.067480492275732
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.054785614195623
$dataBF[17,2] = 1.059921921471272
$dataBF[17,3] = 1.051294384758561
$dataBF[17,4] = 1.067576533879623
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.05434334735113
$dataBF[18,2] = 1.059558968172983
$dataBF[18,3] = 1.050913894465239
$dataBF[18,4] = 1.067147008379679
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.052904820177165
$dataBF[19,2] = 1.058378355601791
$dataBF[19,3] = 1.04967662365388
$dataBF[19,4] = 1.065750483469658
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.051999826402817
$dataBF[20,2] = 1.057635574887797
$dataBF[20,3] = 1.048898492299277
$dataBF[20,4] = 1.064872346161674
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.052479632527245
$dataBF[21,2] = 1.058029383425708
$dataBF[21,3] = 1.049311015063497
$dataBF[21,4] = 1.065337872428373
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.05436745563753
$dataBF[22,2] = 1.059578753265882
$dataBF[22,3] = 1.05093463406882
$dataBF[22,4] = 1.067170420025739
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.056555953691598
$dataBF[23,2] = 1.061374675693944
$dataBF[23,3] = 1.052817908345904
$dataBF[23,4] = 1.069296697072004
$ws.Range("B2:F25").Value = $dataBF

$dataIM = New-Object "object[,]" 24,5
$dataIM[0,0] = 1.049024709399867
$dataIM[0,1] = 1.063291768304708
$dataIM[0,2] = 1.065525881468979
$dataIM[0,3] = 1.057062934024868
$dataIM[0,4] = 1.073690826046026
$dataIM[1,0] = 1.049423142949388
$dataIM[1,1] = 1.064208366506541
$dataIM[1,2] = 1.066378262494624
$dataIM[1,3] = 1.057965346930625
$dataIM[1,4] = 1.074737511571483
$dataIM[2,0] = 1.049679519781438
$dataIM[2,1] = 1.064800548694466
$dataIM[2,2] = 1.066928812330584
$dataIM[2,3] = 1.058548542978494
$dataIM[2,4] = 1.075414216735087
$dataIM[3,0] = 1.049786957525081
$dataIM[3,1] = 1.06504928391776
$dataIM[3,2] = 1.067160026081104
$dataIM[3,3] = 1.058793546390892
$dataIM[3,4] = 1.075698568768255
$dataIM[4,0] = 1.049804976699489
$dataIM[4,1] = 1.065091034910238
$dataIM[4,2] = 1.067198834007942
$dataIM[4,3] = 1.058834673486236
$dataIM[4,4] = 1.075746304870244
$dataIM[5,0] = 1.049680956716911
$dataIM[5,1] = 1.064803873161908
$dataIM[5,2] = 1.066931902751055
$dataIM[5,3] = 1.058551817400908
$dataIM[5,4] = 1.075418016788457
$dataIM[6,0] = 1.049159659752623
$dataIM[6,1] = 1.063601728051215
$dataIM[6,2] = 1.065814155076807
$dataIM[6,3] = 1.057368059783576
$dataIM[6,4] = 1.074044677458174
$dataIM[7,0] = 1.048230032259355
$dataIM[7,1] = 1.06147628931315
$dataIM[7,2] = 1.063836834619052
$dataIM[7,3] = 1.055276505520703
$dataIM[7,4] = 1.071620232053622
$dataIM[8,0] = 1.047602810429161
$dataIM[8,1] = 1.060054447838835
$dataIM[8,2] = 1.062513346522863
$dataIM[8,3] = 1.053878261542364
$dataIM[8,4] = 1.070000834847552
$dataIM[9,0] = 1.047329432996013
$dataIM[9,1] = 1.059437592734008
$dataIM[9,2] = 1.061938989129665
$dataIM[9,3] = 1.053271865792036
$dataIM[9,4] = 1.069298857213029
$dataIM[10,0] = 1.047227619014623
$dataIM[10,1] = 1.059208284581933
$dataIM[10,2] = 1.061725453232275
$dataIM[10,3] = 1.053046479192932
$dataIM[10,4] = 1.069037993800863
$dataIM[11,0] = 1.047249470668729
$dataIM[11,1] = 1.059257480170846
$dataIM[11,2] = 1.061771266216671
$dataIM[11,3] = 1.053094831941229
$dataIM[11,4] = 1.069093955249595
$dataIM[12,0] = 1.04732102252304
$dataIM[12,1] = 1.059418641749396
$dataIM[12,2] = 1.061921342154298
$dataIM[12,3] = 1.0532532382188
$dataIM[12,4] = 1.069277296575545
$dataIM[13,0] = 1.047365072247827
$dataIM[13,1] = 1.05951791470754
$dataIM[13,2] = 1.062013783129237
$dataIM[13,3] = 1.0533508183855
$dataIM[13,4] = 1.069390243557461
$dataIM[14,0] = 1.047620915863902
$dataIM[14,1] = 1.060095361258713
$dataIM[14,2] = 1.062551437633748
$dataIM[14,3] = 1.053918485914416
$dataIM[14,4] = 1.070047406416303
$dataIM[15,0] = 1.047780920824021
$dataIM[15,1] = 1.060457258540783
$dataIM[15,2] = 1.062888350493731
$dataIM[15,3] = 1.054274313869271
$dataIM[15,4] = 1.07045942009947
$dataIM[16,0] = 1.047874076736605
$dataIM[16,1] = 1.060668232590141
$dataIM[16,2] = 1.063084742858169
$dataIM[16,3] = 1.054481771068648
$dataIM[16,4] = 1.070699666769242
$dataIM[17,0] = 1.047905811306888
$dataIM[17,1] = 1.060740149930532
$dataIM[17,2] = 1.063151686739306
$dataIM[17,3] = 1.054552493197637
$dataIM[17,4] = 1.070781572202054
$dataIM[18,0] = 1.047763771633392
$dataIM[18,1] = 1.060418442248568
$dataIM[18,2] = 1.06285221569114
$dataIM[18,3] = 1.054236146348345
$dataIM[18,4] = 1.070415222596184
$dataIM[19,0] = 1.047299959731523
$dataIM[19,1] = 1.059371188686986
$dataIM[19,2] = 1.061877153901821
$dataIM[19,3] = 1.053206595515115
$dataIM[19,4] = 1.069223310391194
$dataIM[20,0] = 1.047006783407008
$dataIM[20,1] = 1.058711692253057
$dataIM[20,2] = 1.061262970534694
$dataIM[20,3] = 1.052558440669093
$dataIM[20,4] = 1.068473225980628
$dataIM[21,0] = 1.047162349891784
$dataIM[21,1] = 1.059061403716321
$dataIM[21,2] = 1.061588667884672
$dataIM[21,3] = 1.05290211965679
$dataIM[21,4] = 1.068870925320781
$dataIM[22,0] = 1.047771521146759
$dataIM[22,1] = 1.060435982014008
$dataIM[22,2] = 1.062868543833719
$dataIM[22,3] = 1.054253392890033
$dataIM[22,4] = 1.070435193773929
$dataIM[23,0] = 1.048471676212392
$dataIM[23,1] = 1.062026618996121
$dataIM[23,2] = 1.064348942159047
$dataIM[23,3] = 1.055817898385703
$dataIM[23,4] = 1.072247547588062
$ws.Range("I2:M25").Value = $dataIM
